# Daily attendance processing - 2026-01-30 09:05:47
# Reorders the "Recorded By" (column G) comma-separated list so that
# "System" (and its variants) appear first, matching the new attendance
# export ordering. Rows whose list includes an admin@admin.com entry, or
# that only contain a single value, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -gt 1 -and ($val -notmatch "admin@admin.com")) {
        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
